# Applies the edit described by the diff:
#  1. Swap the data rows within each year's B/C quarter pair
#     (rows 2<->3, 6<->7, 10<->11, ... 50<->51), across columns A:E.
#  2. Remove columns F (产销率) and G (销售量), which are dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap row pairs (columns A through E) ---
$rowPairs = @()
$r = 2
while ($r -le 50) {
    $rowPairs += , @($r, $r + 1)
    $r += 4
}

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $valsA = @()
    $valsB = @()
    for ($c = 1; $c -le 5; $c++) {
        $valsA += , ($ws.Cells.Item($rowA, $c).Value2)
        $valsB += , ($ws.Cells.Item($rowB, $c).Value2)
    }

    for ($c = 1; $c -le 5; $c++) {
        $oldA = $valsA[$c - 1]
        $oldB = $valsB[$c - 1]
        $newA = $oldB
        $newB = $oldA

        # Assigning "" via .Value deletes the cell outright instead of
        # keeping it as an empty-string cell, so skip no-op empty<->empty
        # swaps to preserve the existing empty inlineStr cell untouched.
        if (-not ($newA -eq "" -and $oldA -eq "")) {
            $ws.Cells.Item($rowA, $c).Value = $newA
        }
        if (-not ($newB -eq "" -and $oldB -eq "")) {
            $ws.Cells.Item($rowB, $c).Value = $newB
        }
    }
}

# --- Step 2: delete columns F and G entirely ---
$ws.Range("F1:G52").Delete() | Out-Null
